$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric stay as text, matching the source data
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.999.54'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '1.743.77'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("D4").Value = '0.9986'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '248.95'
$ws.Range("E5").Value = '  +7.74%  '
$ws.Range("D6").Value = '0.9991'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '0.5157'
$ws.Range("E7").Value = '  -1.45%  '
$ws.Range("D8").Value = '0.2761'
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").Value = '0.06207'
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("D10").Value = '1.742.97'
$ws.Range("E10").Value = '  +0.19%  '
$ws.Range("D11").Value = '0.07219'
$ws.Range("E11").Value = '  +1.24%  '
$ws.Range("D12").Value = '15.18'
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").Value = '0.6515'
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("D14").Value = '4.641'
$ws.Range("E14").Value = '  +2.20%  '
$ws.Range("D15").Value = '77.96'
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("D16").Value = '0.9987'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = '0.9982'
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").Value = '26.030.30'
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").Value = '11.83'
$ws.Range("E19").Value = '  +2.04%  '
$ws.Range("D20").Value = '0.000006827'
$ws.Range("E20").Value = '  +2.54%  '
$ws.Range("D21").Value = '1.967.96'
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").Value = '4.278'
$ws.Range("E22").Value = '  -0.45%  '
$ws.Range("D23").Value = '8.677'
$ws.Range("E23").Value = '  -1.10%  '
$ws.Range("D24").Value = '5.356'
$ws.Range("E24").Value = '  +3.78%  '
$ws.Range("D25").Value = '136.47'
$ws.Range("E25").Value = '  -1.57%  '
$ws.Range("D26").Value = '1.500'
$ws.Range("E26").Value = '  -1.26%  '
$ws.Range("D27").Value = '15.32'
$ws.Range("E27").Value = '  +1.15%  '
$ws.Range("D28").Value = '1.790'
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").Value = '105.48'
$ws.Range("E29").Value = '  +2.26%  '
$ws.Range("D30").Value = '3.954'
$ws.Range("E30").Value = '  +5.36%  '
$ws.Range("D31").Value = '0.08285'
$ws.Range("E31").Value = '  -0.41%  '
$ws.Range("D32").Value = '3.680'
$ws.Range("E32").Value = '  +2.34%  '
$ws.Range("D33").Value = '0.04686'
$ws.Range("E33").Value = '  +3.78%  '
$ws.Range("D34").Value = '2.650'
$ws.Range("E34").Value = '  +1.24%  '
$ws.Range("D35").Value = '1.006'
$ws.Range("E35").Value = '  +2.31%  '
$ws.Range("D36").Value = '0.6245'
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("D37").Value = '2.720'
$ws.Range("E37").Value = '  +1.86%  '
$ws.Range("D38").Value = '0.01603'
$ws.Range("E38").Value = '  +0.52%  '
$ws.Range("D39").Value = '1.942'
$ws.Range("E39").Value = '  +1.59%  '
$ws.Range("D40").Value = '0.9984'
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").Value = '100.89'
$ws.Range("E41").Value = '  +3.25%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.3884'
$ws.Range("E42").Value = '  +0.39%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.7551'
$ws.Range("E43").Value = '  +3.31%  '
$ws.Range("D44").Value = '5.010'
$ws.Range("E44").Value = '  -0.97%  '
$ws.Range("D45").Value = '0.1142'
$ws.Range("E45").Value = '  +1.19%  '
$ws.Range("D46").Value = '6.365'
$ws.Range("E46").Value = '  +1.38%  '
$ws.Range("D47").Value = '55.45'
$ws.Range("E47").Value = '  +3.13%  '
$ws.Range("D48").Value = '0.05222'
$ws.Range("E48").Value = '  -2.21%  '
$ws.Range("D49").Value = '30.63'
$ws.Range("E49").Value = '  +1.25%  '
$ws.Range("D50").Value = '7.665'
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("D51").Value = '0.3446'
$ws.Range("E51").Value = '  +0.50%  '
